$d = $word.ActiveDocument

# The BTec logo's drawing (docPr / cNvPr) is currently named "image1.jpg"
# and the Pearson logo's drawing is currently named "image2.png". The
# edit swaps each drawing's display "name" to the other numeral while
# leaving the "descr" (alt text) and the actual embedded image
# relationships untouched.
$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')
$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')

$d.WordOpenXML = $xml
